$wb = $excel.ActiveWorkbook

# --- Rename the four data sheets to their English short names ---------
$wb.Worksheets.Item("AÑO MODIFICACIÓN").Name = "modified"
$wb.Worksheets.Item("AÑO PUBLICACIÓN").Name = "published"
$wb.Worksheets.Item("TIPO DE OBJETO STIX 2.1 ").Name = "type"
$wb.Worksheets.Item("AÑO CREACION").Name = "created"

# --- Refresh each sheet's embedded chart so its category/value series
#     formulas point at the new sheet name instead of the old one. ------
$wsModified = $wb.Worksheets.Item("modified")
$chartModified = $wsModified.ChartObjects(1).Chart
$serModified = $chartModified.SeriesCollection(1)
$serModified.XValues = "=modified!`$B`$12:`$B`$14"
$serModified.Values  = "=modified!`$D`$12:`$D`$14"

$wsPublished = $wb.Worksheets.Item("published")
$chartPublished = $wsPublished.ChartObjects(1).Chart
$serPublished = $chartPublished.SeriesCollection(1)
$serPublished.XValues = "=published!`$B`$12:`$B`$14"
$serPublished.Values  = "=published!`$D`$12:`$D`$14"

$wsType = $wb.Worksheets.Item("type")
$chartType = $wsType.ChartObjects(1).Chart
$serType = $chartType.SeriesCollection(1)
$serType.XValues = "=type!`$B`$12:`$B`$13"
$serType.Values  = "=type!`$C`$12:`$C`$13"

$wsCreated = $wb.Worksheets.Item("created")
$chartCreated = $wsCreated.ChartObjects(1).Chart
$serCreated = $chartCreated.SeriesCollection(1)
$serCreated.XValues = "=created!`$B`$12:`$B`$15"
$serCreated.Values  = "=created!`$D`$12:`$D`$15"

# --- Move the active/selected tab from "modified" to "created" --------
$wsCreated.Activate()
